# Apply bug-tracker edit: add new bug entry (id 2) to the left-hand "Bugs" table
# (columns A:D) describing the Vertical Position & Horizontal Padding responsiveness
# fix, mirroring the already-resolved first bug's row formatting, then update the
# remaining empty rows (4-10) of that table with the plain bordered style used by
# the right-hand table, and finally move the active selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the cell formatting (styles) into the new cells --------------

# A3 should look like A2 (plain bordered/centered style)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

# B3 should use the plain bordered/centered style (same as A2/F3), not the
# wrap-text style used by B2
$ws.Range("A2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# C3 should look like C2 (yellow "Important" status style)
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# D3 should look like D2 (orange "date fixed" style)
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# Rows 4-10 (columns A:D) get the plain bordered/centered style, same as the
# matching cells in columns F:I
$ws.Range("A2").Copy()
$ws.Range("A4:D10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 2. Fill in the new bug row values -------------------------------------

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Vertical Position & Horizontal Padding ne s'adapte pas à la taille de la vidéo."
$ws.Range("C3").Value = "Important"
$ws.Range("D3").Value = "Corrigé le 16/06/2024"

# --- 3. Update the sheet view (scroll back to A1, move selection) ---------

$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("D5").Select()
